$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the specific cells that were emptied in the edit
$ws.Range("B3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B9").ClearContents()

# Move the active selection to D4
$ws.Range("D4").Select()
